# Apply cryptos price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.951.26"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.641.58"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "`'215.25"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "`'0.5083"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Value = "`'0.2564"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "`'0.06392"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "`'0.07772"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "`'4.299"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "1.646.66"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "`'0.5451"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "0.0₅7838"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "`'64.73"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "25.990.00"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "`'197.59"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").Value = "`'4.429"
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "`'9.946"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "`'6.035"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "`'1.006"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "`'1.869"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").Value = "`'141.18"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "`'0.1143"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "`'6.890"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "`'3.262"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "`'3.185"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "`'1.542"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "`'0.8938"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "`'2.585"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").Value = "1.129.23"
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("D38").Value = "`'0.5503"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "`'0.01553"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "`'1.004"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("B41").Value = "BabyDogeCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D41").Value = "0.0₈131"
$ws.Range("E41").Value = "  +18.97%  "
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").Value = "`'2.549"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "`'0.8162"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "`'100.01"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "1.777.97"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "`'0.4532"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "`'1.005"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "`'54.84"
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "`'0.05079"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "`'1.006"
$ws.Range("E51").Value = "  -0.26%  "
